$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("C14").Value = 0
$ws.Range("C19").Value = 0

$ws.Range("C20").Select()
